$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# [ADD] field offset_footer in mapping
# Add sample "footer" row data illustrating the new offset_footer mapping
# field: the same text appears at different column offsets relative to a
# reference column, with wrap-text enabled so the sample renders on
# multiple lines where the column is narrow.

$ws.Range("C6").Value = "Any footer data"
$ws.Range("C6").WrapText = $true
$ws.Range("E6").WrapText = $true

$ws.Range("E7").Value = "Any footer data"
$ws.Range("E7").WrapText = $true

$ws.Range("F8").Value = "Any footer data"
$ws.Range("F8").WrapText = $true

# Rows with wrapped text in the narrower columns (C, F) grow to fit two
# lines; row 7's text sits in the wider column E and stays single-line.
$ws.Rows.Item(6).RowHeight = 23.85
$ws.Rows.Item(8).RowHeight = 23.85

# Move the active selection to the last edited cell.
[void]$ws.Range("E7").Select()
